# Generate Report for Handback
# Swap the "49a60d07..." and "62542793..." file rows (the handback for
# 49a60d07 is now in and becomes row 2 on every sheet; 62542793 moves to
# row 3), refresh the associated statuses / timestamps / handoff-handback
# file names, and rebuild the hyperlinks whose display text tracks the
# swapped file names.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "49a60d07-fe9b-4c13-ba16-9db3395c7f8e.md"
$ov.Range("B2").Value = "e2e\49a60d07-fe9b-4c13-ba16-9db3395c7f8e.md"
$ov.Range("G2").Value = "2016-08-19 00:46:19"

$ov.Range("A3").Value = "62542793-daaa-495d-a41e-128697f49daa.md"
$ov.Range("B3").Value = "e2e\62542793-daaa-495d-a41e-128697f49daa.md"
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"
$ov.Range("G3").Value = "2016-08-19 00:45:21"

# Rebuild the B2/B3 hyperlinks: same target URLs (rId2 -> 62542793 blob,
# rId3 -> 49a60d07 blob) but the display text now matches the *other*
# file name, mirroring the row swap above.
$ov.Range("B2").Hyperlinks.Delete()
$ov.Range("B3").Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a33d1f88c76e593f0d5f870f41973912325365a5/e2e/62542793-daaa-495d-a41e-128697f49daa.md", "", "", "e2e\49a60d07-fe9b-4c13-ba16-9db3395c7f8e.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fd759b072bc0afe7491fa85fa7b5f86fca7f9a19/e2e/49a60d07-fe9b-4c13-ba16-9db3395c7f8e.md", "", "", "e2e\62542793-daaa-495d-a41e-128697f49daa.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "49a60d07-fe9b-4c13-ba16-9db3395c7f8e.md"
$zh.Range("G2").Value = "49a60d07-fe9b-4c13-ba16-9db3395c7f8e.76fd2a0fcb9275b434388b6bad9e2b7ff4676b21.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-19 00:46:13"
$zh.Range("I2").Value = "49a60d07-fe9b-4c13-ba16-9db3395c7f8e.md"
$zh.Range("J2").Value = "49a60d07-fe9b-4c13-ba16-9db3395c7f8e.76fd2a0fcb9275b434388b6bad9e2b7ff4676b21.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-19 00:46:29"

$zh.Range("A3").Value = "62542793-daaa-495d-a41e-128697f49daa.md"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("G3").Value = "62542793-daaa-495d-a41e-128697f49daa.9782d8192b452272face9fa7282b50c09d296717.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-19 00:45:14"
$zh.Range("I3").Value = "62542793-daaa-495d-a41e-128697f49daa.md"
$zh.Range("J3").Value = "62542793-daaa-495d-a41e-128697f49daa.9782d8192b452272face9fa7282b50c09d296717.zh-cn.xlf"
$zh.Range("P3").ClearContents()

$zh.Range("A2").Hyperlinks.Delete()
$zh.Range("I2").Hyperlinks.Delete()
$zh.Range("A3").Hyperlinks.Delete()
$zh.Range("I3").Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a33d1f88c76e593f0d5f870f41973912325365a5/e2e/62542793-daaa-495d-a41e-128697f49daa.md", "", "", "49a60d07-fe9b-4c13-ba16-9db3395c7f8e.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f1b4dd7e9ebe2ba786c2cdab6aed8648b228f42e/e2e/62542793-daaa-495d-a41e-128697f49daa.md", "", "", "49a60d07-fe9b-4c13-ba16-9db3395c7f8e.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fd759b072bc0afe7491fa85fa7b5f86fca7f9a19/e2e/49a60d07-fe9b-4c13-ba16-9db3395c7f8e.md", "", "", "62542793-daaa-495d-a41e-128697f49daa.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f1b4dd7e9ebe2ba786c2cdab6aed8648b228f42e/e2e/49a60d07-fe9b-4c13-ba16-9db3395c7f8e.md", "", "", "62542793-daaa-495d-a41e-128697f49daa.md") | Out-Null

$zh.Columns.Item(16).ColumnWidth = 12.85

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "49a60d07-fe9b-4c13-ba16-9db3395c7f8e.md"
$de.Range("G2").Value = "49a60d07-fe9b-4c13-ba16-9db3395c7f8e.76fd2a0fcb9275b434388b6bad9e2b7ff4676b21.de-de.xlf"
$de.Range("H2").Value = "2016-08-19 00:46:19"
$de.Range("I2").Value = "49a60d07-fe9b-4c13-ba16-9db3395c7f8e.md"
$de.Range("J2").Value = "49a60d07-fe9b-4c13-ba16-9db3395c7f8e.76fd2a0fcb9275b434388b6bad9e2b7ff4676b21.de-de.xlf"
$de.Range("K2").Value = "2016-08-19 00:46:36"

$de.Range("A3").Value = "62542793-daaa-495d-a41e-128697f49daa.md"
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("G3").Value = "62542793-daaa-495d-a41e-128697f49daa.9782d8192b452272face9fa7282b50c09d296717.de-de.xlf"
$de.Range("H3").Value = "2016-08-19 00:45:21"
$de.Range("I3").Value = "62542793-daaa-495d-a41e-128697f49daa.md"
$de.Range("J3").Value = "62542793-daaa-495d-a41e-128697f49daa.9782d8192b452272face9fa7282b50c09d296717.de-de.xlf"
$de.Range("P3").ClearContents()

$de.Range("A2").Hyperlinks.Delete()
$de.Range("I2").Hyperlinks.Delete()
$de.Range("A3").Hyperlinks.Delete()
$de.Range("I3").Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a33d1f88c76e593f0d5f870f41973912325365a5/e2e/62542793-daaa-495d-a41e-128697f49daa.md", "", "", "49a60d07-fe9b-4c13-ba16-9db3395c7f8e.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/af28aaa2bc22a258483d22ef465041242eaaacb6/e2e/62542793-daaa-495d-a41e-128697f49daa.md", "", "", "49a60d07-fe9b-4c13-ba16-9db3395c7f8e.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fd759b072bc0afe7491fa85fa7b5f86fca7f9a19/e2e/49a60d07-fe9b-4c13-ba16-9db3395c7f8e.md", "", "", "62542793-daaa-495d-a41e-128697f49daa.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/af28aaa2bc22a258483d22ef465041242eaaacb6/e2e/49a60d07-fe9b-4c13-ba16-9db3395c7f8e.md", "", "", "62542793-daaa-495d-a41e-128697f49daa.md") | Out-Null

$de.Columns.Item(16).ColumnWidth = 12.85

Write-Output "done"
